# Updated cryptos list on Wed Jul 17 09:39:30 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ExactText {
    param($cellRange, [string]$text)
    # The source cells are all plain text (coinranking.com scrape
    # dump): prices like "28.80" or "581.03" must stay literal
    # strings, not become numbers (which would normalize away
    # trailing zeros / thousand-dot formatting). A leading
    # apostrophe forces Excel's literal-text entry mode; we then
    # restore the default style so no stray quote-prefix / text
    # number-format remains attached to the cell.
    $cellRange.Value = "'" + $text
    $cellRange.Style = "Normal"
}

$ws.Range("D2").Value = "65.323.19"
$ws.Range("E2").Value = "  +3.61%  "

$ws.Range("D3").Value = "3.491.29"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-ExactText $ws.Range("D5") "581.05"
$ws.Range("E5").Value = "  +2.53%  "

Set-ExactText $ws.Range("D6") "162.34"
$ws.Range("E6").Value = "  +4.53%  "

Set-ExactText $ws.Range("D7") "0.611"
$ws.Range("E7").Value = "  +12.20%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "3.493.40"
$ws.Range("E9").Value = "  +3.02%  "

Set-ExactText $ws.Range("D10") "7.27"
$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("E11").Value = "  +3.67%  "

Set-ExactText $ws.Range("D12") "0.447"
$ws.Range("E12").Value = "  +3.54%  "

$ws.Range("D13").Value = "4.093.16"
$ws.Range("E13").Value = "  +2.95%  "

Set-ExactText $ws.Range("D14") "0.134"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("E15").Value = "  +2.56%  "

Set-ExactText $ws.Range("D16") "28.80"
$ws.Range("E16").Value = "  +6.08%  "

$ws.Range("D17").Value = "65.354.87"
$ws.Range("E17").Value = "  +3.62%  "

$ws.Range("D18").Value = "3.500.37"
$ws.Range("E18").Value = "  +3.36%  "

$ws.Range("E19").Value = "  +3.61%  "

Set-ExactText $ws.Range("D20") "14.42"
$ws.Range("E20").Value = "  +2.50%  "

Set-ExactText $ws.Range("D21") "385.20"
$ws.Range("E21").Value = "  +2.16%  "

$ws.Range("E22").Value = "  +2.64%  "

$ws.Range("E23").Value = "  +4.59%  "

Set-ExactText $ws.Range("D24") "72.74"
$ws.Range("E24").Value = "  +2.15%  "

Set-ExactText $ws.Range("D25") "0.998"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  +3.19%  "

$ws.Range("E27").Value = "  +7.62%  "

$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("E29").Value = "  -0.07%  "

Set-ExactText $ws.Range("D30") "1.53"
$ws.Range("E30").Value = "  +13.66%  "

Set-ExactText $ws.Range("D31") "6.19"
$ws.Range("E31").Value = "  +1.73%  "

$ws.Range("E32").Value = "  +3.47%  "

Set-ExactText $ws.Range("D33") "23.73"
$ws.Range("E33").Value = "  +2.49%  "

Set-ExactText $ws.Range("D34") "7.19"
$ws.Range("E34").Value = "  +6.21%  "

$ws.Range("E35").Value = "  +12.30%  "

Set-ExactText $ws.Range("D36") "162.62"
$ws.Range("E36").Value = "  +1.92%  "

Set-ExactText $ws.Range("D37") "1.93"
$ws.Range("E37").Value = "  +6.03%  "

$ws.Range("D38").Value = "3.034.16"
$ws.Range("E38").Value = "  +2.69%  "

$ws.Range("E39").Value = "  +3.57%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-ExactText $ws.Range("D40") "6.90"
$ws.Range("E40").Value = "  +8.75%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-ExactText $ws.Range("D41") "27.08"
$ws.Range("E41").Value = "  +0.46%  "

Set-ExactText $ws.Range("D42") "4.61"
$ws.Range("E42").Value = "  +6.82%  "

$ws.Range("E43").Value = "  +1.24%  "

Set-ExactText $ws.Range("D44") "43.01"
$ws.Range("E44").Value = "  +3.40%  "

Set-ExactText $ws.Range("D45") "0.782"
$ws.Range("E45").Value = "  +3.27%  "

Set-ExactText $ws.Range("D46") "25.91"
$ws.Range("E46").Value = "  +11.43%  "

Set-ExactText $ws.Range("D48") "320.67"
$ws.Range("E48").Value = "  +11.03%  "

Set-ExactText $ws.Range("D49") "6.76"
$ws.Range("E49").Value = "  +6.49%  "

Set-ExactText $ws.Range("D50") "0.883"
$ws.Range("E50").Value = "  +6.59%  "

$ws.Range("E51").Value = "  +6.66%  "

